# Update the "fixed" date text (shown in the Date placeholders of the
# slide master, every slide layout, and the notes master) from 5/6/21 to
# 6/15/21, and bump the version label "v1.0" -> "v2.0" on the 4 version
# badge shapes on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "5/6/21"
$newDate = "6/15/21"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Version badge shapes on slide 1: "v1.0" -> "v2.0"
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "v1.0") {
            $tr.Text = "v2.0"
        }
    }
}
